# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
#
# Several match rows in the "Croatia 1NL" sheet got re-paired: the betting
# odds / result data (columns B, E:AD - i.e. everything except the row
# index column A, the Div column C and the Date column D) that belonged to
# one fixture were moved to a different row. This script re-creates that
# re-shuffle by reading the affected rows into memory first and then
# writing them back out in their new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and E..AD (everything that "travels" with a fixture, excluding
# A = row index, C = Div, D = Date, which stay put because the swapped
# fixtures share the same matchday/date).
$cols = @(2,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30)

function Get-RowData($r) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Cells.Item($r, $c).Value2
    }
    return $data
}

function Set-RowData($r, $data) {
    foreach ($c in $cols) {
        $val = $data[$c]
        if ($null -eq $val) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value2 = $val
        }
    }
}

# --- Rows 9, 10, 11: 3-way rotation (row9 <- row10 <- row11 <- row9) ---
$r9  = Get-RowData 9
$r10 = Get-RowData 10
$r11 = Get-RowData 11

Set-RowData 9  $r10
Set-RowData 10 $r11
Set-RowData 11 $r9

# --- Simple 2-row swaps ---
$swapPairs = @(
    @(52, 54),
    @(70, 71),
    @(131, 132),
    @(142, 144),
    @(203, 204)
)

foreach ($pair in $swapPairs) {
    $a = $pair[0]
    $b = $pair[1]
    $rowA = Get-RowData $a
    $rowB = Get-RowData $b
    Set-RowData $a $rowB
    Set-RowData $b $rowA
}
